$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Append a new registration row (row 18) below the existing data.
$newRow = 18

$ws.Cells.Item($newRow, 1).Value = "Parceira"
$ws.Cells.Item($newRow, 2).Value = "Gildean Nunes Mendes "
$ws.Cells.Item($newRow, 4).Value = "Usimig"
$ws.Cells.Item($newRow, 5).Value = "B1 - Substituir Caçamba Recuperadora Tipo Ponte"

# Column F ("Data") holds plain text like "2025-12-15" in this sheet, not a
# real date. Force text interpretation so Excel doesn't auto-convert it to a
# date serial, then drop the temporary Text number format so the cell keeps
# the workbook's default (unstyled) appearance, matching the other rows.
$ws.Cells.Item($newRow, 6).NumberFormat = "@"
$ws.Cells.Item($newRow, 6).Value = "2025-12-15"
$ws.Cells.Item($newRow, 6).ClearFormats()

$ws.Cells.Item($newRow, 7).Value = "ADM (09-16h)"
$ws.Cells.Item($newRow, 8).Value = "Turno C"
